# Add Setups and Results
# Replace the placeholder DOE setup/result values with the computed
# full-factorial setups (POX/C, C/A, POX/M) and their derived results.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Low/high levels for each factor (C/A and POX/M are the POX/C level
# scaled down by 1e3 and 1e5 respectively).
$POXC_low  = 80.40036015459947
$POXC_high = 119.5996398454005
$CA_low    = 0.08040036015459946
$CA_high   = 0.1195996398454005
$POXM_low  = 0.0008040036015459946
$POXM_high = 0.001195996398454005

# Full-factorial design matrix (rows 2-9), columns B:D = POX/C, C/A, POX/M
$rows = @(
    @{ Row = 2; POXC = $POXC_low;  CA = $CA_low;  POXM = $POXM_low  },
    @{ Row = 3; POXC = $POXC_high; CA = $CA_low;  POXM = $POXM_low  },
    @{ Row = 4; POXC = $POXC_low;  CA = $CA_high; POXM = $POXM_low  },
    @{ Row = 5; POXC = $POXC_high; CA = $CA_high; POXM = $POXM_low  },
    @{ Row = 6; POXC = $POXC_low;  CA = $CA_low;  POXM = $POXM_high },
    @{ Row = 7; POXC = $POXC_high; CA = $CA_low;  POXM = $POXM_high },
    @{ Row = 8; POXC = $POXC_low;  CA = $CA_high; POXM = $POXM_high },
    @{ Row = 9; POXC = $POXC_high; CA = $CA_high; POXM = $POXM_high }
)

foreach ($r in $rows) {
    $ws.Cells.Item($r.Row, 2).Value = $r.POXC
    $ws.Cells.Item($r.Row, 3).Value = $r.CA
    $ws.Cells.Item($r.Row, 4).Value = $r.POXM
}
